$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 4 and 5 swap their per-record data (the underlying observation that
# used to be on row 4 is now on row 5, and vice versa). Columns P, S, T, U,
# V, W, AD, AE, AG, AT, AW, AY are identical between the two rows already,
# so only the differing columns need to be exchanged:
# A, B, D, E, F, G, H, Q, R, Y, AA, AX
# ---------------------------------------------------------------------------

# Y and AA hold dates stored as plain text (e.g. "2025-08-29"). Force the
# cells to a text number format first so Excel does not reinterpret the
# swapped text as a date serial number when it is written back.
foreach ($addr in @("Y4", "AA4", "Y5", "AA5")) {
    $ws.Range($addr).NumberFormat = "@"
}

$cols4_5 = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Y", "AA", "AX")
foreach ($col in $cols4_5) {
    $rng4 = $ws.Range("$col`4")
    $rng5 = $ws.Range("$col`5")
    $tmp = $rng4.Value2
    $rng4.Value2 = $rng5.Value2
    $rng5.Value2 = $tmp
}

# ---------------------------------------------------------------------------
# Rows 13 and 14 swap their per-record data the same way. Here the Y, AA and
# AX columns already hold identical values on both rows, so they do not need
# to be touched.
# A, B, D, E, F, G, H, Q, R
# ---------------------------------------------------------------------------
$cols13_14 = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")
foreach ($col in $cols13_14) {
    $rng13 = $ws.Range("$col`13")
    $rng14 = $ws.Range("$col`14")
    $tmp = $rng13.Value2
    $rng13.Value2 = $rng14.Value2
    $rng14.Value2 = $tmp
}
